# Weekly refresh of the Ají (Hortaliza) price log for
# "Vega Monumental Concepción": insert a new, more recent price record at
# row 10 and push the existing 31 data rows (old rows 10-40) down by one
# (new rows 11-41), growing the used range from A1:R40 to A1:R41.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 10; Excel shifts rows 10..40 down to 11..41
# and extends the sheet dimension to A1:R41 automatically.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly record.
$ws.Cells.Item(10, 1).Value  = 11
$ws.Cells.Item(10, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(10, 3).Value  = "Bíobío"
$ws.Cells.Item(10, 4).Value  = 44453
$ws.Cells.Item(10, 5).Value  = 8
$ws.Cells.Item(10, 6).Value  = 100112021
$ws.Cells.Item(10, 7).Value  = "Ají"
$ws.Cells.Item(10, 8).Value  = "Americana (o)"
$ws.Cells.Item(10, 9).Value  = "Primera"
$ws.Cells.Item(10, 10).Value = 50
$ws.Cells.Item(10, 11).Value = 35000
$ws.Cells.Item(10, 12).Value = 36000
$ws.Cells.Item(10, 13).Value = 35400
$ws.Cells.Item(10, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(10, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(10, 16).Value = 2950
$ws.Cells.Item(10, 17).Value = 12
$ws.Cells.Item(10, 18).Value = "Hortaliza"
